$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the new "notes" column (D) for the header + first three data rows,
# matching the formatting of column C in the same row.
$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("D2").Value = "notes"

$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").Value = "First note"

$ws.Range("C4").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("D4").Value = "First note"

$ws.Range("C5").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("D5").Value = "First note"

# Add a new data row (row 6), duplicating the "shashwat" record with a second note
$ws.Range("A4:D4").Copy()
$ws.Range("A6:D6").PasteSpecial(-4122)

$ws.Range("C6").Copy()
$ws.Range("D6").PasteSpecial(-4122)

$ws.Range("A6").Value = "shashwat"
$ws.Range("B6").Value = 123
$ws.Range("C6").Value = "Kalkaji"
$ws.Range("D6").Value = "Second note"
